$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had columns I (Baricenter + Greedy Crossings) and J
# (BranchAndReduce Duration [ns]); both are dropped, and the old column K
# (BranchAndReduce Crossings) shifts left to become the new column I.
$ws.Range("I1:J14").EntireColumn.Delete()

# Rename the remaining headers for the new column meanings.
$ws.Range("D1").Value = "Barycenter Duration [ns]"
$ws.Range("E1").Value = "Barycenter Crossings"
$ws.Range("F1").Value = "Median Duration [ns]"
$ws.Range("G1").Value = "Median Crossings"
$ws.Range("H1").Value = "BranchAndReduce Duration [ns]"
# I1 already reads "BranchAndReduce Crossings" after the column delete above.

# Row 2: matching_4_4.gr
$ws.Range("B2").Value = "'2535"
$ws.Range("D2").Value = "'532"
$ws.Range("E2").Value = "'0"
$ws.Range("F2").Value = "'416"
$ws.Range("H2").Value = "'7725"

# Row 3: cycle_8_sorted.gr
$ws.Range("B3").Value = "'2536"
$ws.Range("D3").Value = "'249"
$ws.Range("F3").Value = "'449"
$ws.Range("H3").Value = "'18996"

# Row 4: tree_6_10.gr
$ws.Range("B4").Value = "'383294106"
$ws.Range("D4").Value = "'696"
$ws.Range("E4").Value = "'13"
$ws.Range("F4").Value = "'908"
$ws.Range("G4").Value = "'13"
$ws.Range("H4").Value = "'23267"

# Row 5: cycle_8_shuffled.gr
$ws.Range("B5").Value = "'3102"
$ws.Range("D5").Value = "'184"
$ws.Range("F5").Value = "'334"
$ws.Range("H5").Value = "'9256"

# Row 6: complete_4_5.gr
$ws.Range("B6").Value = "'15397"
$ws.Range("D6").Value = "'311"
$ws.Range("F6").Value = "'483"
$ws.Range("H6").Value = "'3497"

# Row 7: path_9_shuffled.gr
$ws.Range("B7").Value = "'2704"
$ws.Range("D7").Value = "'154"
$ws.Range("E7").Value = "'6"
$ws.Range("F7").Value = "'225"
$ws.Range("H7").Value = "'8480"

# Row 8: ladder_4_4_sorted.gr
$ws.Range("B8").Value = "'4053"
$ws.Range("D8").Value = "'243"
$ws.Range("E8").Value = "'3"
$ws.Range("F8").Value = "'413"
$ws.Range("G8").Value = "'3"
$ws.Range("H8").Value = "'7139"

# Row 9: ladder_4_4_shuffled.gr
$ws.Range("B9").Value = "'3582"
$ws.Range("D9").Value = "'245"
$ws.Range("F9").Value = "'422"
$ws.Range("H9").Value = "'6033"

# Row 10: path_9_sorted.gr
$ws.Range("B10").Value = "'2781"
$ws.Range("D10").Value = "'138"
$ws.Range("E10").Value = "'0"
$ws.Range("F10").Value = "'282"
$ws.Range("H10").Value = "'5123"

# Row 11: website_20.gr
$ws.Range("B11").Value = "'372812331"
$ws.Range("D11").Value = "'576"
$ws.Range("E11").Value = "'17"
$ws.Range("F11").Value = "'772"
$ws.Range("H11").Value = "'25228"

# Row 12: star_6.gr
$ws.Range("B12").Value = "'27331"
$ws.Range("D12").Value = "'217"
$ws.Range("E12").Value = "'0"
$ws.Range("F12").Value = "'293"
$ws.Range("H12").Value = "'7101"

# Row 13: plane_5_6.gr
$ws.Range("B13").Value = "'59038"
$ws.Range("D13").Value = "'373"
$ws.Range("E13").Value = "'0"
$ws.Range("F13").Value = "'501"
$ws.Range("G13").Value = "'0"
$ws.Range("H13").Value = "'6823"

# Row 14: grid_9_shuffled.gr
$ws.Range("B14").Value = "'8429"
$ws.Range("D14").Value = "'274"
$ws.Range("E14").Value = "'17"
$ws.Range("F14").Value = "'455"
$ws.Range("G14").Value = "'21"
$ws.Range("H14").Value = "'9808"
